# Apply the "excel demo and some other classes added" edit to LoginList.xlsx
# Target sheet: QaTeam3 (3rd sheet / sheet3.xml) gains a new D:I block of
# "school program" data, becomes the active tab, and its selection moves to I15.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("QaTeam3")

# --- Header row (D1:I1) -----------------------------------------------------
# These header cells share the same style as A1:C1 (12pt font -> cellXf s="2").
$ws.Range("D1").Value = "okul"
$ws.Range("E1").Value = "degree"
$ws.Range("F1").Value = "study"
$ws.Range("G1").Value = "fromdate"
$ws.Range("H1").Value = "Todate"
$ws.Range("I1").Value = "Program Desc"
$ws.Range("D1:I1").Font.Size = 12

# --- Data, entered column-by-column (matches the shared-string insert order) -
$schools = @("gazi", "otdü", "bilkent", "ege", "itü", "ktü", "yildiz", "istanbul", "ankara")
for ($i = 0; $i -lt $schools.Length; $i++) {
    $ws.Cells.Item($i + 2, 4).Value = $schools[$i]
}

$degrees = @("good", "nice", "very good", "good", "nice", "very good", "good", "nice", "very good")
for ($i = 0; $i -lt $degrees.Length; $i++) {
    $ws.Cells.Item($i + 2, 5).Value = $degrees[$i]
}

$studies = @("finance", "medical", "engineer", "finance", "medical", "engineer", "finance", "medical", "engineer")
for ($i = 0; $i -lt $studies.Length; $i++) {
    $ws.Cells.Item($i + 2, 6).Value = $studies[$i]
}

$fromDates = @(11122021, 11122020, 11122019, 11122021, 11122020, 11122019, 11122021, 11122020, 11122019)
for ($i = 0; $i -lt $fromDates.Length; $i++) {
    $ws.Cells.Item($i + 2, 7).Value = $fromDates[$i]
}

$toDates = @(12122022, 12122022, 12122022, 12122022, 12122022, 12122022, 12122022, 12122022, 12122022)
for ($i = 0; $i -lt $toDates.Length; $i++) {
    $ws.Cells.Item($i + 2, 8).Value = $toDates[$i]
}

$progs = @("Nice Program", "Excellent Program", "Not Bad program", "Nice Program", "Excellent Program", "Not Bad program", "Nice Program", "Excellent Program", "Not Bad program")
for ($i = 0; $i -lt $progs.Length; $i++) {
    $ws.Cells.Item($i + 2, 9).Value = $progs[$i]
}

# --- Column widths for the new columns --------------------------------------
$ws.Columns.Item(7).ColumnWidth = 10
$ws.Columns.Item(8).ColumnWidth = 8.67
$ws.Columns.Item(9).ColumnWidth = 14.67

# --- Make QaTeam3 the active sheet/tab and move its selection to I15 --------
$ws.Activate() | Out-Null
$ws.Range("I15").Select() | Out-Null
